$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1631
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""

$ws.Range("H36").Value = 1631
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = ""

$ws.Range("H70").Value = 1987.125
$ws.Range("I70").Value = 1566.6666
$ws.Range("J70").Value = 2239.4
$ws.Range("K70").Value = 4699.9998
$ws.Range("L70").Value = 6718.200000000001
$ws.Range("M70").Value = -4429.9998
$ws.Range("N70").Value = -7258.200000000001

$ws.Range("H73").Value = 1987.125
$ws.Range("I73").Value = 1566.6666
$ws.Range("J73").Value = 2239.4
$ws.Range("K73").Value = 4699.9998
$ws.Range("L73").Value = 6718.200000000001
$ws.Range("M73").Value = -3763.9998
$ws.Range("N73").Value = -8590.200000000001

$ws.Range("H132").Value = 4033.3704
$ws.Range("I132").Value = 2815.125
$ws.Range("J132").Value = 5805.364
$ws.Range("K132").Value = 8445.375
$ws.Range("L132").Value = 17416.092
$ws.Range("M132").Value = -5915.375
$ws.Range("N132").Value = -22476.092

$ws.Range("H138").Value = 3734.6316
$ws.Range("I138").Value = 3965.5715
$ws.Range("K138").Value = 11896.7145
$ws.Range("M138").Value = -6756.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2965.25
$ws.Range("I2").Value = 616.3333
$ws.Range("K2").Value = 616.3333
$ws.Range("M2").Value = -503.3333

$ws.Range("H50").Value = 39316.668
$ws.Range("I50").Value = 38000
$ws.Range("J50").Value = 39975
$ws.Range("K50").Value = 38000
$ws.Range("L50").Value = 39975
$ws.Range("M50").Value = -37286
$ws.Range("N50").Value = -41403

$ws.Range("H88").Value = 4288.5386
$ws.Range("I88").Value = 2965.6667
$ws.Range("J88").Value = 5422.4287
$ws.Range("K88").Value = 2965.6667
$ws.Range("L88").Value = 5422.4287
$ws.Range("M88").Value = -2559.6667
$ws.Range("N88").Value = -6234.4287

$ws.Range("H91").Value = 4288.5386
$ws.Range("I91").Value = 2965.6667
$ws.Range("J91").Value = 5422.4287
$ws.Range("K91").Value = 2965.6667
$ws.Range("L91").Value = 5422.4287
$ws.Range("M91").Value = -1561.6667
$ws.Range("N91").Value = -8230.4287

$ws.Range("H97").Value = 1666.6666
$ws.Range("I97").Value = 666.6667
$ws.Range("J97").Value = 2666.6667
$ws.Range("K97").Value = 666.6667
$ws.Range("L97").Value = 2666.6667
$ws.Range("M97").Value = -170.6667
$ws.Range("N97").Value = -3658.6667

$ws.Range("H116").Value = 2965.25
$ws.Range("I116").Value = 616.3333
$ws.Range("K116").Value = 616.3333
$ws.Range("M116").Value = 1677.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2965.25
$ws.Range("I3").Value = 616.3333
$ws.Range("K3").Value = 616.3333
$ws.Range("M3").Value = -502.3333

$ws.Range("H20").Value = 4463
$ws.Range("I20").Value = 2670.75
$ws.Range("J20").Value = 8047.5
$ws.Range("K20").Value = 2670.75
$ws.Range("L20").Value = 8047.5
$ws.Range("M20").Value = -2423.75
$ws.Range("N20").Value = -8541.5

$ws.Range("H56").Value = 39999
$ws.Range("J56").Value = 39999
$ws.Range("L56").Value = 39999
$ws.Range("N56").Value = -41477

$ws.Range("H94").Value = 3169.7273
$ws.Range("I94").Value = 1644.5
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1644.5
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -1193.5
$ws.Range("N94").Value = -5902

$ws.Range("H99").Value = 2286.125
$ws.Range("I99").Value = 1912.7142
$ws.Range("K99").Value = 1912.7142
$ws.Range("M99").Value = -414.7141999999999

$ws.Range("H134").Value = 4010.923
$ws.Range("J134").Value = 4605.625
$ws.Range("L134").Value = 13816.875
$ws.Range("N134").Value = -18886.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 559.44446
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

$ws.Range("H55").Value = 7090.909
$ws.Range("I55").Value = 6800
$ws.Range("K55").Value = 6800
$ws.Range("M55").Value = -6485

$ws.Range("H122").Value = 468
$ws.Range("I122").Value = 468
$ws.Range("K122").Value = 1404
$ws.Range("M122").Value = 1046

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 3161
$ws.Range("I33").Value = 500
$ws.Range("J33").Value = 4491.5
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 26949
$ws.Range("M33").Value = -2717
$ws.Range("N33").Value = -27515

$ws.Range("H50").Value = 340.54544
$ws.Range("I50").Value = 224
$ws.Range("J50").Value = 1506
$ws.Range("K50").Value = 672
$ws.Range("L50").Value = 4518
$ws.Range("M50").Value = -191
$ws.Range("N50").Value = -5480

$ws.Range("H53").Value = 340.54544
$ws.Range("I53").Value = 224
$ws.Range("J53").Value = 1506
$ws.Range("K53").Value = 672
$ws.Range("L53").Value = 4518
$ws.Range("M53").Value = -191
$ws.Range("N53").Value = -5480

$ws.Range("H60").Value = 4060
$ws.Range("I60").Value = 122.5
$ws.Range("K60").Value = 367.5
$ws.Range("M60").Value = -116.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3836
$ws.Range("I70").Value = 3254
$ws.Range("K70").Value = 3254
$ws.Range("M70").Value = -2984

$ws.Range("H73").Value = 3836
$ws.Range("I73").Value = 3254
$ws.Range("K73").Value = 3254
$ws.Range("M73").Value = -2318

$ws.Range("H113").Value = 6000
$ws.Range("J113").Value = 6000
$ws.Range("L113").Value = 6000
$ws.Range("N113").Value = -10340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 22000
$ws.Range("J38").Value = 22000
$ws.Range("L38").Value = 22000
$ws.Range("N38").Value = -22820

$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4540

$ws.Range("H46").Value = 2664.7334
$ws.Range("I46").Value = 2068.1875
$ws.Range("K46").Value = 2068.1875
$ws.Range("M46").Value = -1880.1875

$ws.Range("H53").Value = 13511.5
$ws.Range("I53").Value = 13515.333
$ws.Range("K53").Value = 13515.333
$ws.Range("M53").Value = -12997.333

$ws.Range("H55").Value = 623.875
$ws.Range("I55").Value = 58.6
$ws.Range("K55").Value = 58.6
$ws.Range("M55").Value = 114.4

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = ""
$ws.Range("N56").Value = ""

$ws.Range("H69").Value = 66666
$ws.Range("J69").Value = 66666
$ws.Range("L69").Value = 66666
$ws.Range("N69").Value = -68288

$ws.Range("H72").Value = 66666
$ws.Range("J72").Value = 66666
$ws.Range("L72").Value = 199998
$ws.Range("N72").Value = -208110

$ws.Range("H82").Value = 2411.0557
$ws.Range("I82").Value = 940.8
$ws.Range("J82").Value = 2976.5386
$ws.Range("K82").Value = 940.8
$ws.Range("L82").Value = 2976.5386
$ws.Range("M82").Value = -579.8
$ws.Range("N82").Value = -3698.5386

$ws.Range("H85").Value = 2411.0557
$ws.Range("I85").Value = 940.8
$ws.Range("J85").Value = 2976.5386
$ws.Range("K85").Value = 940.8
$ws.Range("L85").Value = 2976.5386
$ws.Range("M85").Value = 307.2
$ws.Range("N85").Value = -5472.5386

$ws.Range("H101").Value = 26794.6
$ws.Range("J101").Value = 26794.6
$ws.Range("L101").Value = 26794.6
$ws.Range("N101").Value = -33284.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""

$ws.Range("H136").Value = 549.10345
$ws.Range("I136").Value = 552.7778
$ws.Range("J136").Value = 499.5
$ws.Range("K136").Value = 1658.3334
$ws.Range("L136").Value = 1498.5
$ws.Range("M136").Value = 891.6666
$ws.Range("N136").Value = -6598.5
